$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Andre B.-Eletrônica analóg. e de potên"
$ws.Range("D2").Value = "[André Guimarães-CAD, -]"
$ws.Range("F2").Value = "[-, João Paulo-Sistemas digitais, João Paulo-Sistemas digitais]"

# Row 3
$ws.Range("B3").Value = "Nilton Maia-M.T.R"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "[-, -, Pedro Bispo-Acionamentos Elétricos, -]"
$ws.Range("F3").Value = "[Sandro-Programação de Computadores, -, -, -]"

# Row 4
$ws.Range("B4").Value = "Nilton Maia-M.T.R"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "[-, -, Pedro Bispo-Acionamentos Elétricos, -]"
$ws.Range("F4").Value = "[Sandro-Programação de Computadores, -, -, -]"

# Row 6
$ws.Range("B6").Value = "Sandro-Circuitos elétrico"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "[-, -, Pedro Bispo-Acionamentos Elétricos, -]"
$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "[Sandro-Programação de Computadores, -, -, -]"

# Row 7
$ws.Range("B7").Value = "Sandro-Circuitos elétrico"
$ws.Range("C7").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "[Sandro-Programação de Computadores, -, -, -]"

# Row 8
$ws.Range("B8").Value = "Andre B.-Eletrônica analóg. e de potên"
$ws.Range("D8").Value = "[-, André Guimarães-CAD]"
$ws.Range("F8").Value = "[Pedro Bispo-Acionamentos Elétricos, -, -, -]"
